$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.584.28"
$ws.Range("E2").Value = "  +0.90%  "

$ws.Range("D3").Value = "1.873.61"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.28"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4728"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2905"
$ws.Range("E8").Value = "  +1.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06474"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.10"
$ws.Range("E10").Value = "  +4.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07705"
$ws.Range("E11").Value = "  -0.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7380"
$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.20"
$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("D14").Value = "1.873.28"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.153"
$ws.Range("E15").Value = "  +0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.24"
$ws.Range("E16").Value = "  -0.47%  "

$ws.Range("D17").Value = "30.651.70"
$ws.Range("E17").Value = "  +1.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.30"
$ws.Range("E18").Value = "  -0.58%  "

$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007505"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("D21").Value = "2.119.69"
$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.259"
$ws.Range("E23").Value = "  +0.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.182"
$ws.Range("E24").Value = "  +0.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.199"
$ws.Range("E25").Value = "  -0.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.90"
$ws.Range("E26").Value = "  -0.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.72"
$ws.Range("E27").Value = "  -0.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.905"
$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("E29").Value = "  +1.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.344"
$ws.Range("E30").Value = "  -2.58%  "

$ws.Range("E31").Value = "  -0.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.269"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.094"
$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04797"
$ws.Range("E34").Value = "  +0.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.121"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6945"
$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.748"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.230"
$ws.Range("E40").Value = "  -3.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.04"
$ws.Range("E41").Value = "  +4.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.967"
$ws.Range("E42").Value = "  +2.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4168"
$ws.Range("E43").Value = "  +1.20%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8356"
$ws.Range("E45").Value = "  -0.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.14"
$ws.Range("E46").Value = "  -0.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.378"
$ws.Range("E47").Value = "  -0.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.41"
$ws.Range("E48").Value = "  +0.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.991"
$ws.Range("E49").Value = "  -1.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "915.70"
$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05654"
$ws.Range("E51").Value = "  +1.52%  "
